$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: wrap the existing remarks text (style changes to wrap text)
$ws.Range("H16").WrapText = $true

# Row 17: fill in the second sprint-review table row
$ws.Range("B17").Value = "Erstellung der Home-site"
$ws.Range("C17").Value = 1.5
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("H17").WrapText = $true
$ws.Range("H17").Value = "Ich fand es gut das wir dass MockUp erstellt haben sodass jeder weiß wie die Seiten ungefähr ausschauen sollen. Das wir jetzt auch ein Kanban board haben ist auch gut so weiß jeder was gemacht werden muss und man hat generell eine klare übersicht über das Projekt . Beim erstellen der Homeseite gab es wenig probleme weil ich schon das wissen dafür hatte aus dem letzten Semester. Wir hätten das Mockup und Das Kanbanboard erstellen sollen. Ich finde wir sollten in zunkunft mehere Gruppen meetings machen."

# Row heights grow because of the long wrapped remarks text
$ws.Rows("16:17").RowHeight = 194.25

# Scroll/selection state left behind after entering the data
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("K17").Select()
